$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewGroupStruct")

# Insert a new row at position 34 (shifts old rows 34-38 down to 35-39),
# mirroring the "TypeValue" field being added to the TimeZoneGroupDto struct
# documented in this second copy of the table.
$ws.Rows.Item(34).Insert() | Out-Null
$ws.Rows.Item(34).RowHeight = 18.75

# New row 34: TypeValue / string / string
$ws.Range("B34").Value = "TypeValue"
$ws.Range("C34").Value = "string"
$ws.Range("D34").Value = "string"

# Fill in new example values in column D for the existing fields.
$ws.Range("D30").Value = '""'
$ws.Range("D31").Value = "America"
$ws.Range("D32").Value = "America"
$ws.Range("D33").Value = "americaTimeZones"

# Old row 34 (IanaVariableName) is now row 35; give it an example value too.
$ws.Range("D35").Value = "americaTimeZones"

# Widen column C a bit to fit the longer example values.
$ws.Columns.Item(3).ColumnWidth = 31.88

# Touch the page setup (portrait orientation) for this sheet.
$ws.PageSetup.Orientation = 1

# Make this the active sheet/tab and select D37, matching the saved view.
$ws.Activate() | Out-Null
$ws.Range("D37").Select() | Out-Null
